$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "46.919.41"
$ws.Cells.Item(2, 5).Value = "  +5.11%  "

$ws.Cells.Item(3, 4).Value = "2.349.87"
$ws.Cells.Item(3, 5).Value = "  +4.55%  "

$ws.Cells.Item(4, 5).Value = "  -0.58%  "

Set-TextValue $ws.Cells.Item(5, 4) "306.87"
$ws.Cells.Item(5, 5).Value = "  -0.01%  "

Set-TextValue $ws.Cells.Item(6, 4) "98.27"
$ws.Cells.Item(6, 5).Value = "  +3.45%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.578"
$ws.Cells.Item(7, 5).Value = "  +1.27%  "

Set-TextValue $ws.Cells.Item(8, 4) "1.00"
$ws.Cells.Item(8, 5).Value = "  -0.55%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.538"
$ws.Cells.Item(9, 5).Value = "  +4.19%  "

Set-TextValue $ws.Cells.Item(10, 4) "35.77"
$ws.Cells.Item(10, 5).Value = "  +2.36%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.0806"
$ws.Cells.Item(11, 5).Value = "  +0.54%  "

$ws.Cells.Item(12, 5).Value = "  +3.37%  "

$ws.Cells.Item(13, 5).Value = "  -0.38%  "

$ws.Cells.Item(14, 4).Value = "2.706.94"
$ws.Cells.Item(14, 5).Value = "  +4.46%  "

$ws.Cells.Item(15, 4).Value = "2.349.62"
$ws.Cells.Item(15, 5).Value = "  +4.76%  "

Set-TextValue $ws.Cells.Item(16, 4) "14.24"
$ws.Cells.Item(16, 5).Value = "  +5.01%  "

Set-TextValue $ws.Cells.Item(17, 4) "0.833"
$ws.Cells.Item(17, 5).Value = "  +0.04%  "

$ws.Cells.Item(18, 4).Value = "46.785.74"
$ws.Cells.Item(18, 5).Value = "  +5.31%  "

Set-TextValue $ws.Cells.Item(19, 4) "13.77"
$ws.Cells.Item(19, 5).Value = "  +17.48%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0950"
$ws.Cells.Item(20, 5).Value = "  +1.45%  "

Set-TextValue $ws.Cells.Item(21, 4) "6.19"
$ws.Cells.Item(21, 5).Value = "  -0.10%  "

Set-TextValue $ws.Cells.Item(22, 4) "66.88"
$ws.Cells.Item(22, 5).Value = "  +2.33%  "

Set-TextValue $ws.Cells.Item(23, 4) "245.37"
$ws.Cells.Item(23, 5).Value = "  +3.29%  "

Set-TextValue $ws.Cells.Item(24, 4) "2.99"
$ws.Cells.Item(24, 5).Value = "  +1.32%  "

Set-TextValue $ws.Cells.Item(25, 4) "2.00"
$ws.Cells.Item(25, 5).Value = "  +1.08%  "

Set-TextValue $ws.Cells.Item(26, 4) "0.992"
$ws.Cells.Item(26, 5).Value = "  -1.11%  "

Set-TextValue $ws.Cells.Item(27, 4) "41.90"
$ws.Cells.Item(27, 5).Value = "  +13.30%  "

$ws.Cells.Item(28, 5).Value = "  -0.68%  "

$ws.Cells.Item(29, 5).Value = "  +1.47%  "

Set-TextValue $ws.Cells.Item(30, 4) "20.18"
$ws.Cells.Item(30, 5).Value = "  +0.91%  "

Set-TextValue $ws.Cells.Item(31, 4) "5.78"
$ws.Cells.Item(31, 5).Value = "  -1.99%  "

Set-TextValue $ws.Cells.Item(32, 4) "152.55"
$ws.Cells.Item(32, 5).Value = "  +3.14%  "

$ws.Cells.Item(33, 5).Value = "  +3.97%  "

$ws.Cells.Item(34, 5).Value = "  +0.34%  "

$ws.Cells.Item(35, 5).Value = "  -1.01%  "

$ws.Cells.Item(37, 5).Value = "  +0.58%  "

$ws.Cells.Item(38, 5).Value = "  -1.00%  "

Set-TextValue $ws.Cells.Item(39, 4) "4.05"
$ws.Cells.Item(39, 5).Value = "  +7.03%  "

$ws.Cells.Item(40, 5).Value = "  +5.68%  "

Set-TextValue $ws.Cells.Item(41, 4) "3.43"
$ws.Cells.Item(41, 5).Value = "  +2.55%  "

Set-TextValue $ws.Cells.Item(42, 4) "13.82"
$ws.Cells.Item(42, 5).Value = "  -9.17%  "

$ws.Cells.Item(43, 4).Value = "1.878.57"
$ws.Cells.Item(43, 5).Value = "  +3.89%  "

Set-TextValue $ws.Cells.Item(44, 4) "0.999"
$ws.Cells.Item(44, 5).Value = "  -0.81%  "

Set-TextValue $ws.Cells.Item(45, 4) "1.99"
$ws.Cells.Item(45, 5).Value = "  +13.00%  "

Set-TextValue $ws.Cells.Item(46, 4) "0.197"
$ws.Cells.Item(46, 5).Value = "  +5.20%  "

Set-TextValue $ws.Cells.Item(47, 4) "74.47"
$ws.Cells.Item(47, 5).Value = "  +8.01%  "

Set-TextValue $ws.Cells.Item(48, 4) "81.01"
$ws.Cells.Item(48, 5).Value = "  -1.18%  "

Set-TextValue $ws.Cells.Item(49, 4) "4.94"

Set-TextValue $ws.Cells.Item(50, 4) "99.01"
$ws.Cells.Item(50, 5).Value = "  +0.62%  "

Set-TextValue $ws.Cells.Item(51, 4) "55.53"
$ws.Cells.Item(51, 5).Value = "  +2.74%  "
